$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.919.09'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.293.83'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.83'
$ws.Range('E5').Value = '  +18.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.73'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.08'
$ws.Range('E10').Value = '  +5.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0948'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.97'
$ws.Range('E12').Value = '  +12.90%  '
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.89'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.635.87'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.853'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.291.83'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.769.57'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.76'
$ws.Range('E20').Value = '  +9.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.35'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.44'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.90'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +9.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.62'
$ws.Range('E25').Value = '  +6.04%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.65'
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.87'
$ws.Range('E28').Value = '  +8.48%  '
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.13'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0933'
$ws.Range('E32').Value = '  +4.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.56'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.69'
$ws.Range('E34').Value = '  +4.63%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.66'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0363'
$ws.Range('E37').Value = '  +3.25%  '
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.78'
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.244'
$ws.Range('E40').Value = '  +3.52%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.01'
$ws.Range('E41').Value = '  +14.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.38'
$ws.Range('E42').Value = '  +3.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.62'
$ws.Range('E43').Value = '  +10.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.30'
$ws.Range('E44').Value = '  +21.28%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.40'
$ws.Range('E46').Value = '  +3.98%  '
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0996'
$ws.Range('E48').Value = '  -2.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.68'
$ws.Range('E49').Value = '  +4.26%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.24'
$ws.Range('E50').Value = '  +3.06%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.469'
$ws.Range('E51').Value = '  +6.99%  '
